$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.000119846596356663
$ws.Range("C2").Value = 0.000119846596356663
$ws.Range("D2").Value = 0.983041706615532
$ws.Range("E2").Value = 0.999161073825503
$ws.Range("F2").Value = 0.000958772770853308
$ws.Range("G2").Value = 0.00743048897411314
$ws.Range("H2").Value = 0.00035953978906999
$ws.Range("I2").Value = 0.000299616490891659
$ws.Range("J2").Value = 0.991670661553212
$ws.Range("K2").Value = 0.999220997123682
$ws.Range("L2").Value = 0.000179769894534995
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0.000479386385426654
$ws.Range("O2").Value = 0.999101150527325
$ws.Range("P2").Value = 0.000119846596356663
$ws.Range("Q2").Value = 0.999580536912752
$ws.Range("R2").Value = 0.0000599232981783317
$ws.Range("S2").Value = 0.998382070949185
$ws.Range("T2").Value = 0.000239693192713327
$ws.Range("U2").Value = 0.000119846596356663
$ws.Range("V2").Value = 0.947267497603068
$ws.Range("W2").Value = 0.000179769894534995
$ws.Range("X2").Value = 0.000119846596356663
$ws.Range("B3").Value = 0.999520613614573
$ws.Range("C3").Value = 0.999820230105465
$ws.Range("D3").Value = 0.00185762224352828
$ws.Range("E3").Value = 0.000419463087248322
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.0000599232981783317
$ws.Range("J3").Value = 0.000179769894534995
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0.0000599232981783317
$ws.Range("M3").Value = 0.0225311601150527
$ws.Range("N3").Value = 0.0000599232981783317
$ws.Range("O3").Value = 0.000299616490891659
$ws.Range("P3").Value = 0.890280441035475
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0.000299616490891659
$ws.Range("S3").Value = 0.00119846596356663
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0.999460690316395
$ws.Range("V3").Value = 0.000119846596356663
$ws.Range("W3").Value = 0.0000599232981783317
$ws.Range("X3").Value = 0.0000599232981783317
$ws.Range("B4").Value = 0.0000599232981783317
$ws.Range("C4").Value = 0.0000599232981783317
$ws.Range("D4").Value = 0.00383509108341323
$ws.Range("E4").Value = 0.000179769894534995
$ws.Range("F4").Value = 0.998981303930968
$ws.Range("G4").Value = 0.992150047938639
$ws.Range("H4").Value = 0.99964046021093
$ws.Range("I4").Value = 0.99964046021093
$ws.Range("J4").Value = 0.00725071907957814
$ws.Range("K4").Value = 0.000779002876318313
$ws.Range("L4").Value = 0.999760306807287
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.999460690316395
$ws.Range("O4").Value = 0.000479386385426654
$ws.Range("P4").Value = 0.0000599232981783317
$ws.Range("Q4").Value = 0.000179769894534995
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0.000239693192713327
$ws.Range("T4").Value = 0.999700383509108
$ws.Range("U4").Value = 0.0000599232981783317
$ws.Range("V4").Value = 0.0506951102588687
$ws.Range("W4").Value = 0.999580536912752
$ws.Range("X4").Value = 0.999700383509108
$ws.Range("B5").Value = 0.000299616490891659
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.0100671140939597
$ws.Range("E5").Value = 0.000119846596356663
$ws.Range("F5").Value = 0.0000599232981783317
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0.000179769894534995
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0.976929530201342
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0.000119846596356663
$ws.Range("P5").Value = 0.105644774688399
$ws.Range("Q5").Value = 0.000239693192713327
$ws.Range("R5").Value = 0.99964046021093
$ws.Range("S5").Value = 0.000179769894534995
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0.000239693192713327
$ws.Range("V5").Value = 0.00119846596356663
$ws.Range("W5").Value = 0.0000599232981783317
$ws.Range("X5").Value = 0.0000599232981783317

$wb.Save()
